# Apply cryptos list update (Tue May  2 20:13:39 UTC 2023, GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.797.15'
$ws.Range("E2").Value = '  +2.94%  '

# Row 3
$ws.Range("D3").Value = '1.879.77'
$ws.Range("E3").Value = '  +3.19%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.51%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.48'
$ws.Range("E5").Value = '  -0.48%  '

# Row 6
$ws.Range("E6").Value = '  +0.32%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4665'
$ws.Range("E7").Value = '  +0.60%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3946'
$ws.Range("E8").Value = '  +2.70%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07921'
$ws.Range("E9").Value = '  +1.19%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9767'
$ws.Range("E10").Value = '  +1.90%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.39'
$ws.Range("E11").Value = '  +2.30%  '

# Row 12
$ws.Range("D12").Value = '1.873.81'
$ws.Range("E12").Value = '  -0.34%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.758'
$ws.Range("E13").Value = '  +1.94%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.012'
$ws.Range("E14").Value = '  +2.26%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06981'
$ws.Range("E15").Value = '  +1.76%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.64'
$ws.Range("E16").Value = '  +2.42%  '

# Row 17
$ws.Range("E17").Value = '  +0.43%  '

# Row 18
$ws.Range("E18").Value = '  +1.84%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.02'
$ws.Range("E19").Value = '  +1.70%  '

# Row 20
$ws.Range("E20").Value = '  +0.14%  '

# Row 21
$ws.Range("D21").Value = '28.813.79'
$ws.Range("E21").Value = '  +2.85%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.363'
$ws.Range("E22").Value = '  +0.80%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.11'
$ws.Range("E23").Value = '  +1.61%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.122'
$ws.Range("E24").Value = '  +1.20%  '

# Row 25
$ws.Range("D25").Value = '2.134.98'
$ws.Range("E25").Value = '  +3.09%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.73'
$ws.Range("E26").Value = '  +1.41%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.44'
$ws.Range("E27").Value = '  +1.33%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.756'
$ws.Range("E28").Value = '  +0.95%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.009'
$ws.Range("E29").Value = '  +2.10%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.92'
$ws.Range("E30").Value = '  +2.95%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09387'
$ws.Range("E31").Value = '  +1.36%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9427'
$ws.Range("E32").Value = '  +0.41%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.324'
$ws.Range("E33").Value = '  +0.92%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.353'
$ws.Range("E34").Value = '  +3.31%  '

# Row 35
$ws.Range("E35").Value = '  -1.94%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05915'
$ws.Range("E36").Value = '  -0.45%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02125'
$ws.Range("E37").Value = '  -0.99%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.149'
$ws.Range("E38").Value = '  -0.23%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.937'
$ws.Range("E39").Value = '  +5.02%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5707'
$ws.Range("E40").Value = '  +2.09%  '

# Row 41
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.989'
$ws.Range("E41").Value = '  +0.92%  '

# Row 42
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1793'
$ws.Range("E42").Value = '  +1.38%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07250'
$ws.Range("E43").Value = '  +3.52%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.82'
$ws.Range("E44").Value = '  +2.18%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5349'
$ws.Range("E45").Value = '  +1.82%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.156'
$ws.Range("E46").Value = '  -6.51%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.125'
$ws.Range("E47").Value = '  -5.08%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.853'
$ws.Range("E48").Value = '  +1.54%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '114.33'
$ws.Range("E49").Value = '  +1.72%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.373'
$ws.Range("E50").Value = '  +3.03%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
$ws.Range("E51").Value = '  +0.42%  '
